# Updated symbol list on Mon Jan  9 07:30:35 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures on Sheet1.
# Values are stored as literal text (matching the source data's inline-string
# cells), so each cell is briefly put into Text format before the write and
# restored to the default "Normal" style afterwards to avoid leaving any
# formatting residue behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.61%"
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.41%"
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.820"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.07%"
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06245"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.48%"
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.924"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.57%"
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.266"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.09%"
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8796"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.31%"
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.82%"
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1451"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.46%"
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05253"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.53%"
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.75%"
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03084"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.46%"
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09053"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.03%"
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.47%"
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006262"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.55%"
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005911"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.16%"
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.452"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.20%"
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "6.50%"
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3131"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.78%"
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1294"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.28%"
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.845"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.81%"
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04312"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.45%"
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001178"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.78%"
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.13%"
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.08%"
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001689"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.02%"
$ws.Range("E27").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04037"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.05%"
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006703"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "62.13%"
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1154"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.80%"
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002133"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.53%"
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.26%"
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005090"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.39%"
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.07%"
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.377"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "839.17%"
$ws.Range("E47").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("E50").Style = "Normal"
